# Project Updated at: 2024-07-25 21-45-31
#
# Applies the texts.xlsx update: refreshes stale "Wildcard
# Characters"/"Wildcard Ranges" columns on the Typography sheet,
# fixes stale IT/ES/FR placeholders on Translation row 637, and
# populates newly-used Translation rows 654-663.

$wb = $excel.ActiveWorkbook
$typo = $wb.Worksheets.Item("Typography")
$trans = $wb.Worksheets.Item("Translation")

# --- Typography sheet: rows 6, 12, 24, 25, 26 had a stale/blank ---
# --- "Wildcard Characters" (G) and "Wildcard Ranges" (I) value;  ---
# --- refresh them to match the value already used on rows 4-5.  ---
$wildcardChars = '!”"#*%&()''$+-@_, .:;?/~±×÷•º`´{}©£€^®¥_=[]¡¢|\¿><'
$wildcardRanges = 'a-z,A-Z,0-9'

foreach ($r in 6,12,24,25,26) {
    $typo.Cells.Item($r, 7).Value = $wildcardChars
    $typo.Cells.Item($r, 9).Value = $wildcardRanges
}

# Excel auto-detects a bare "0" as a number, but the source data
# keeps it as text (matching column F on the same rows), so force
# a temporary Text format for the write and reset the style back
# to Normal afterwards (leaves no visible/format change behind).

# --- Translation sheet: row 637 G/H/I (IT/ES/FR) stale "-" -> "0" ---
foreach ($c in 7,8,9) {
    $cell = $trans.Cells.Item(637, $c)
    $cell.NumberFormat = "@"
    $cell.Value = "0"
    $cell.Style = "Normal"
}

# --- Translation sheet: populate newly-used rows 654-663 ---
# (previously blank placeholder rows in the sheet)

# row 654
$trans.Cells.Item(654, 2).Value = 'SingleUseId4054'
$trans.Cells.Item(654, 3).Value = 'Verdana_20_itaic'
$trans.Cells.Item(654, 4).Value = 'Center'
$trans.Cells.Item(654, 5).Value = 'LTR'
$trans.Cells.Item(654, 6).Value = '<value>'
$trans.Cells.Item(654, 7).Value = '<value>'
$trans.Cells.Item(654, 8).Value = '<value>'
$trans.Cells.Item(654, 9).Value = '<value>'

# row 655
$trans.Cells.Item(655, 2).Value = 'SingleUseId4055'
$trans.Cells.Item(655, 3).Value = 'Verdana_20_itaic'
$trans.Cells.Item(655, 4).Value = 'Left'
$trans.Cells.Item(655, 5).Value = 'LTR'
$trans.Cells.Item(655, 6).Value = 'Finalizado!'
$trans.Cells.Item(655, 7).Value = 'Finalizado!'
$trans.Cells.Item(655, 8).Value = 'Finalizado!'
$trans.Cells.Item(655, 9).Value = 'Finalizado!'

# row 656
$trans.Cells.Item(656, 2).Value = 'SingleUseId4056'
$trans.Cells.Item(656, 3).Value = 'Default'
$trans.Cells.Item(656, 4).Value = 'Center'
$trans.Cells.Item(656, 5).Value = 'LTR'
$trans.Cells.Item(656, 6).Value = '<value>'
$trans.Cells.Item(656, 7).Value = '<value>'
$trans.Cells.Item(656, 8).Value = '<value>'
$trans.Cells.Item(656, 9).Value = '<value>'

# row 657
$trans.Cells.Item(657, 2).Value = 'SingleUseId4057'
$trans.Cells.Item(657, 3).Value = 'Default'
$trans.Cells.Item(657, 4).Value = 'Left'
$trans.Cells.Item(657, 5).Value = 'LTR'
$cell = $trans.Cells.Item(657, 6)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(657, 7)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(657, 8)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(657, 9)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"

# row 658
$trans.Cells.Item(658, 2).Value = 'SingleUseId4058'
$trans.Cells.Item(658, 3).Value = 'Default'
$trans.Cells.Item(658, 4).Value = 'Center'
$trans.Cells.Item(658, 5).Value = 'LTR'
$trans.Cells.Item(658, 6).Value = '<value>'
$trans.Cells.Item(658, 7).Value = '<value>'
$trans.Cells.Item(658, 8).Value = '<value>'
$trans.Cells.Item(658, 9).Value = '<value>'

# row 659
$trans.Cells.Item(659, 2).Value = 'SingleUseId4059'
$trans.Cells.Item(659, 3).Value = 'Default'
$trans.Cells.Item(659, 4).Value = 'Left'
$trans.Cells.Item(659, 5).Value = 'LTR'
$cell = $trans.Cells.Item(659, 6)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(659, 7)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(659, 8)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(659, 9)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"

# row 660
$trans.Cells.Item(660, 2).Value = 'SingleUseId4060'
$trans.Cells.Item(660, 3).Value = 'Default'
$trans.Cells.Item(660, 4).Value = 'Center'
$trans.Cells.Item(660, 5).Value = 'LTR'
$trans.Cells.Item(660, 6).Value = '<value>'
$trans.Cells.Item(660, 7).Value = '<value>'
$trans.Cells.Item(660, 8).Value = '<value>'
$trans.Cells.Item(660, 9).Value = '<value>'

# row 661
$trans.Cells.Item(661, 2).Value = 'SingleUseId4061'
$trans.Cells.Item(661, 3).Value = 'Default'
$trans.Cells.Item(661, 4).Value = 'Left'
$trans.Cells.Item(661, 5).Value = 'LTR'
$cell = $trans.Cells.Item(661, 6)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(661, 7)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(661, 8)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(661, 9)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"

# row 662
$trans.Cells.Item(662, 2).Value = 'SingleUseId4062'
$trans.Cells.Item(662, 3).Value = 'Default'
$trans.Cells.Item(662, 4).Value = 'Center'
$trans.Cells.Item(662, 5).Value = 'LTR'
$trans.Cells.Item(662, 6).Value = '<value>'
$trans.Cells.Item(662, 7).Value = '<value>'
$trans.Cells.Item(662, 8).Value = '<value>'
$trans.Cells.Item(662, 9).Value = '<value>'

# row 663
$trans.Cells.Item(663, 2).Value = 'SingleUseId4063'
$trans.Cells.Item(663, 3).Value = 'Default'
$trans.Cells.Item(663, 4).Value = 'Left'
$trans.Cells.Item(663, 5).Value = 'LTR'
$cell = $trans.Cells.Item(663, 6)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(663, 7)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(663, 8)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"
$cell = $trans.Cells.Item(663, 9)
$cell.NumberFormat = "@"
$cell.Value = '0'
$cell.Style = "Normal"

Write-Output "edit applied"
